# Add a new timeline row to the end of the first (and only) table,
# describing the diamond-square algorithm prototype work session.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()

# Col 1: Date of Session
$newRow.Cells.Item(1).Range.Text = "14/11/2021"

# Col 2: Time Spent
$newRow.Cells.Item(2).Range.Text = "7 Hours"

# Col 3: Development Segment (needs its own shading colour, distinct
# from the "System Design" blue inherited from the row above)
$newRow.Cells.Item(3).Range.Text = "World Generation – Objective 1"
$newRow.Cells.Item(3).Shading.BackgroundPatternColor = 4765306

# Col 4: Notes
$newRow.Cells.Item(4).Range.Text = "Attempted to create a prototype for the diamond-square algorithm to be imported into the main system. This unfortunately was largely a failure but can be built from in the future."
